# Applies the cryptos.xlsx diff (crypto price/volume refresh + Stellar/Monero row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.481.81"
$ws.Range("E2").Value = "  -2.65%  "
$ws.Range("D3").Value = "3.258.81"
$ws.Range("E3").Value = "  -5.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.02"
$ws.Range("E5").Value = "  -2.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.11"
$ws.Range("E6").Value = "  -10.85%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.249.91"
$ws.Range("E8").Value = "  -5.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.542"
$ws.Range("E9").Value = "  -8.76%  "
$ws.Range("E10").Value = "  -11.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.70"
$ws.Range("E11").Value = "  -4.44%  "
$ws.Range("E12").Value = "  -10.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.30"
$ws.Range("E13").Value = "  -13.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000245"
$ws.Range("E14").Value = "  -8.88%  "
$ws.Range("D15").Value = "3.778.89"
$ws.Range("E15").Value = "  -5.43%  "
$ws.Range("D16").Value = "67.506.50"
$ws.Range("E16").Value = "  -2.66%  "
$ws.Range("D17").Value = "3.258.30"
$ws.Range("E17").Value = "  -5.38%  "
$ws.Range("E18").Value = "  -5.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "527.07"
$ws.Range("E19").Value = "  -9.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.09"
$ws.Range("E20").Value = "  -12.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.93"
$ws.Range("E21").Value = "  -12.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.753"
$ws.Range("E22").Value = "  -10.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.83"
$ws.Range("E23").Value = "  -12.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.36"
$ws.Range("E24").Value = "  -11.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.46"
$ws.Range("E25").Value = "  -11.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.21"
$ws.Range("E27").Value = "  -11.08%  "
$ws.Range("E28").Value = "  -11.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.98"
$ws.Range("E29").Value = "  -7.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.90"
$ws.Range("E30").Value = "  -11.70%  "
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("E32").Value = "  -4.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.59"
$ws.Range("E33").Value = "  -15.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.67"
$ws.Range("E34").Value = "  -13.38%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "511.69"
$ws.Range("E36").Value = "  -11.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0440"
$ws.Range("E37").Value = "  -6.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "53.09"
$ws.Range("E38").Value = "  -5.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0851"
$ws.Range("E39").Value = "  -10.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.89"
$ws.Range("E40").Value = "  -15.37%  "
$ws.Range("E41").Value = "  -11.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.76"
$ws.Range("E42").Value = "  -11.94%  "
$ws.Range("D43").Value = "2.929.92"
$ws.Range("E43").Value = "  -9.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.265"
$ws.Range("E44").Value = "  -10.27%  "
$ws.Range("D45").Value = "0.0₃0586"
$ws.Range("E45").Value = "  -14.33%  "
$ws.Range("E46").Value = "  -8.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.47"
$ws.Range("E47").Value = "  -14.73%  "
$ws.Range("E49").Value = "  -16.79%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "123.78"
$ws.Range("E50").Value = "  -7.62%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.113"
$ws.Range("E51").Value = "  -10.34%  "
